$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 5668.263
$ws.Range("I62").Value = 4209.5454
$ws.Range("J62").Value = 7674
$ws.Range("K62").Value = 4209.5454
$ws.Range("L62").Value = 7674
$ws.Range("M62").Value = -3585.5454
$ws.Range("N62").Value = -8922
# Row 64
$ws.Range("H64").Value = 3548.8333
$ws.Range("I64").Value = 2971.75
$ws.Range("J64").Value = 3837.375
$ws.Range("K64").Value = 2971.75
$ws.Range("L64").Value = 3837.375
$ws.Range("M64").Value = -2723.75
$ws.Range("N64").Value = -4333.375
# Row 65
$ws.Range("H65").Value = 5668.263
$ws.Range("I65").Value = 4209.5454
$ws.Range("J65").Value = 7674
$ws.Range("K65").Value = 21047.727
$ws.Range("L65").Value = 38370
$ws.Range("M65").Value = -17927.727
$ws.Range("N65").Value = -44610
# Row 67
$ws.Range("H67").Value = 3548.8333
$ws.Range("I67").Value = 2971.75
$ws.Range("J67").Value = 3837.375
$ws.Range("K67").Value = 2971.75
$ws.Range("L67").Value = 3837.375
$ws.Range("M67").Value = -2113.75
$ws.Range("N67").Value = -5553.375
# Row 86
$ws.Range("H86").Value = 5423.0435
$ws.Range("I86").Value = 1152.5333
$ws.Range("J86").Value = 13430.25
$ws.Range("K86").Value = 1152.5333
$ws.Range("L86").Value = 13430.25
$ws.Range("M86").Value = -29.53330000000005
$ws.Range("N86").Value = -15676.25
# Row 89
$ws.Range("H89").Value = 5423.0435
$ws.Range("I89").Value = 1152.5333
$ws.Range("J89").Value = 13430.25
$ws.Range("K89").Value = 5762.6665
$ws.Range("L89").Value = 67151.25
$ws.Range("M89").Value = -146.6665000000003
$ws.Range("N89").Value = -78383.25
# Row 98
$ws.Range("H98").Value = 707.56525
$ws.Range("I98").Value = 467.42856
$ws.Range("K98").Value = 467.42856
$ws.Range("M98").Value = 1030.57144
# Row 101
$ws.Range("H101").Value = 423.33334
$ws.Range("I101").Value = 300
$ws.Range("J101").Value = 485
$ws.Range("K101").Value = 900
$ws.Range("L101").Value = 1455
$ws.Range("M101").Value = 722
$ws.Range("N101").Value = -4699
# Row 111
$ws.Range("H111").Value = 4366.2144
$ws.Range("I111").Value = 5125.222
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 15375.666
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -12308.666
$ws.Range("N111").Value = -15134
# Row 122
$ws.Range("H122").Value = 707.56525
$ws.Range("I122").Value = 467.42856
$ws.Range("K122").Value = 1402.28568
$ws.Range("M122").Value = 1047.71432
# Row 132
$ws.Range("H132").Value = 1901.8793
$ws.Range("I132").Value = 2119.2126
$ws.Range("J132").Value = 973.2727
$ws.Range("K132").Value = 6357.6378
$ws.Range("L132").Value = 2919.8181
$ws.Range("M132").Value = -3827.6378
$ws.Range("N132").Value = -7979.8181

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4696.375
$ws.Range("I32").Value = 3961.589
$ws.Range("K32").Value = 3961.589
$ws.Range("M32").Value = -3674.589
# Row 43
$ws.Range("H43").Value = 11999.5
$ws.Range("J43").Value = 11999.5
$ws.Range("L43").Value = 11999.5
$ws.Range("N43").Value = -12625.5
# Row 102
$ws.Range("H102").Value = 1213.375
$ws.Range("I102").Value = 1084.5
$ws.Range("K102").Value = 1084.5
$ws.Range("M102").Value = 537.5
# Row 122
$ws.Range("H122").Value = 1066.2373
$ws.Range("I122").Value = 927.84
$ws.Range("J122").Value = 1835.1111
$ws.Range("K122").Value = 2783.52
$ws.Range("L122").Value = 5505.3333
$ws.Range("M122").Value = -333.52
$ws.Range("N122").Value = -10405.3333
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1774.75
$ws.Range("I99").Value = 1866.3334
$ws.Range("K99").Value = 1866.3334
$ws.Range("M99").Value = -368.3334
# Row 107
$ws.Range("H107").Value = 908.3125
$ws.Range("I107").Value = 518.5
$ws.Range("J107").Value = 2077.75
$ws.Range("K107").Value = 518.5
$ws.Range("L107").Value = 2077.75
$ws.Range("M107").Value = 1401.5
$ws.Range("N107").Value = -5917.75

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 99
$ws.Range("H99").Value = 3453.0476
$ws.Range("I99").Value = 2857.1428
$ws.Range("J99").Value = 4644.857
$ws.Range("K99").Value = 2857.1428
$ws.Range("L99").Value = 4644.857
$ws.Range("M99").Value = -1359.1428
$ws.Range("N99").Value = -7640.857
# Row 126
$ws.Range("H126").Value = 3453.0476
$ws.Range("I126").Value = 2857.1428
$ws.Range("J126").Value = 4644.857
$ws.Range("K126").Value = 8571.428400000001
$ws.Range("L126").Value = 13934.571
$ws.Range("M126").Value = -6101.428400000001
$ws.Range("N126").Value = -18874.571
# Row 132
$ws.Range("H132").Value = 2716.7334
$ws.Range("I132").Value = 1923.2727
$ws.Range("J132").Value = 4898.75
$ws.Range("K132").Value = 5769.8181
$ws.Range("L132").Value = 14696.25
$ws.Range("M132").Value = -3239.8181
$ws.Range("N132").Value = -19756.25

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 5735.6665
$ws.Range("I107").Value = 6423.9375
$ws.Range("J107").Value = 229.5
$ws.Range("K107").Value = 19271.8125
$ws.Range("L107").Value = 688.5
$ws.Range("M107").Value = -17351.8125
$ws.Range("N107").Value = -4528.5
# Row 131
$ws.Range("H131").Value = 115741.02
$ws.Range("I131").Value = 840
$ws.Range("J131").Value = 118444.58
$ws.Range("K131").Value = 2520
$ws.Range("L131").Value = 355333.74
$ws.Range("M131").Value = 2520
$ws.Range("N131").Value = -365413.74

$ws = $wb.Worksheets.Item("GSM")
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 113
$ws.Range("H113").Value = 17542.2
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4064.2144
$ws.Range("I7").Value = 3911.111
$ws.Range("J7").Value = 4339.8
$ws.Range("K7").Value = 3911.111
$ws.Range("L7").Value = 4339.8
$ws.Range("M7").Value = -3799.111
$ws.Range("N7").Value = -4563.8
# Row 61
$ws.Range("H61").Value = 3474.611
$ws.Range("I61").Value = 1702.8667
$ws.Range("J61").Value = 12333.333
$ws.Range("K61").Value = 1702.8667
$ws.Range("L61").Value = 12333.333
$ws.Range("M61").Value = -1500.8667
$ws.Range("N61").Value = -12737.333
# Row 93
$ws.Range("H93").Value = 1075.75
$ws.Range("I93").Value = 849.5
$ws.Range("J93").Value = 1302
$ws.Range("K93").Value = 849.5
$ws.Range("L93").Value = 1302
$ws.Range("M93").Value = 398.5
$ws.Range("N93").Value = -3798
# Row 113
$ws.Range("H113").Value = 3474.611
$ws.Range("I113").Value = 1702.8667
$ws.Range("J113").Value = 12333.333
$ws.Range("K113").Value = 1702.8667
$ws.Range("L113").Value = 12333.333
$ws.Range("M113").Value = 467.1333
$ws.Range("N113").Value = -16673.333
# Row 126
$ws.Range("H126").Value = 4064.2144
$ws.Range("I126").Value = 3911.111
$ws.Range("J126").Value = 4339.8
$ws.Range("K126").Value = 11733.333
$ws.Range("L126").Value = 13019.4
$ws.Range("M126").Value = -9263.332999999999
$ws.Range("N126").Value = -17959.4

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 142857890
$ws.Range("I107").Value = 250000210
$ws.Range("J107").Value = 1467
$ws.Range("K107").Value = 750000630
$ws.Range("L107").Value = 4401
$ws.Range("M107").Value = -749998710
$ws.Range("N107").Value = -8241
# Row 132
$ws.Range("H132").Value = 1165.122
$ws.Range("I132").Value = 743.44446
$ws.Range("J132").Value = 1978.3572
$ws.Range("K132").Value = 2230.33338
$ws.Range("L132").Value = 5935.071599999999
$ws.Range("M132").Value = 299.66662
$ws.Range("N132").Value = -10995.0716
# Row 136
$ws.Range("H136").Value = 21068716
$ws.Range("I136").Value = 29494116
$ws.Range("J136").Value = 5213.5713
$ws.Range("K136").Value = 88482348
$ws.Range("L136").Value = 15640.7139
$ws.Range("M136").Value = -88479798
$ws.Range("N136").Value = -20740.7139
